$values = @(
    "audio0",
    "soundfiles/NU1109_0021_talker01_SNR0.wav",
    "soundfiles/NU1109_0040_talker01_SNR0.wav",
    "soundfiles/NU1109_0049_talker01_SNR0.wav",
    "soundfiles/NU1109_0059_talker01_SNR0.wav",
    "soundfiles/NU1109_0087_talker01_SNR0.wav",
    "soundfiles/NU1109_0097_talker01_SNR0.wav",
    "soundfiles/NU1109_0101_talker01_SNR0.wav",
    "soundfiles/NU1109_0108_talker01_SNR0.wav",
    "soundfiles/NU1109_0121_talker01_SNR0.wav",
    "soundfiles/NU1109_0125_talker01_SNR0.wav",
    "soundfiles/NU1109_0146_talker01_SNR0.wav",
    "soundfiles/NU1109_0148_talker01_SNR0.wav",
    "soundfiles/NU1109_0153_talker01_SNR0.wav",
    "soundfiles/NU1109_0155_talker01_SNR0.wav",
    "soundfiles/NU1109_0163_talker01_SNR0.wav",
    "soundfiles/NU1109_0170_talker01_SNR0.wav",
    "soundfiles/NU1109_0179_talker01_SNR0.wav",
    "soundfiles/NU1109_0180_talker01_SNR0.wav",
    "soundfiles/NU1109_0183_talker01_SNR0.wav",
    "soundfiles/NU1109_0198_talker01_SNR0.wav",
    "soundfiles/NU1109_0212_talker01_SNR0.wav",
    "soundfiles/NU1109_0216_talker01_SNR0.wav",
    "soundfiles/NU1109_0219_talker01_SNR0.wav",
    "soundfiles/NU1109_0221_talker01_SNR0.wav",
    "soundfiles/NU1109_0227_talker01_SNR0.wav",
    "soundfiles/NU1109_0243_talker01_SNR0.wav",
    "soundfiles/NU1109_0244_talker01_SNR0.wav",
    "soundfiles/NU1109_0245_talker01_SNR0.wav",
    "soundfiles/NU1109_0248_talker01_SNR0.wav",
    "soundfiles/NU1109_0249_talker01_SNR0.wav",
    "soundfiles/NU1109_0255_talker01_SNR0.wav",
    "soundfiles/NU1109_0261_talker01_SNR0.wav",
    "soundfiles/NU1109_0274_talker01_SNR0.wav",
    "soundfiles/NU1109_0285_talker01_SNR0.wav",
    "soundfiles/NU1109_0287_talker01_SNR0.wav",
    "soundfiles/NU1109_0288_talker01_SNR0.wav",
    "soundfiles/NU1109_0291_talker01_SNR0.wav",
    "soundfiles/NU1109_0292_talker01_SNR0.wav",
    "soundfiles/NU1109_0293_talker01_SNR0.wav",
    "soundfiles/NU1109_0297_talker01_SNR0.wav",
    "soundfiles/NU1109_0306_talker01_SNR0.wav",
    "soundfiles/NU1109_0310_talker01_SNR0.wav",
    "soundfiles/NU1109_0319_talker01_SNR0.wav",
    "soundfiles/NU1109_0322_talker01_SNR0.wav",
    "soundfiles/NU1109_0326_talker01_SNR0.wav",
    "soundfiles/NU1109_0329_talker01_SNR0.wav",
    "soundfiles/NU1109_0341_talker01_SNR0.wav",
    "soundfiles/NU1109_0344_talker01_SNR0.wav",
    "soundfiles/NU1109_0345_talker01_SNR0.wav",
    "soundfiles/NU1109_0359_talker01_SNR0.wav",
    "soundfiles/NU1109_0361_talker01_SNR0.wav",
    "soundfiles/NU1109_0365_talker01_SNR0.wav",
    "soundfiles/NU1109_0374_talker01_SNR0.wav",
    "soundfiles/NU1109_0376_talker01_SNR0.wav",
    "soundfiles/NU1109_0377_talker01_SNR0.wav",
    "soundfiles/NU1109_0386_talker01_SNR0.wav",
    "soundfiles/NU1109_0387_talker01_SNR0.wav",
    "soundfiles/NU1109_0389_talker01_SNR0.wav",
    "soundfiles/NU1109_0391_talker01_SNR0.wav",
    "soundfiles/NU1109_0392_talker01_SNR0.wav"
)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $cell = $ws.Cells.Item($row, 3)
    $cell.Value = $values[$i]
}

# Match formatting of column B for the new column C cells
$ws.Range("C1:C61").NumberFormat = $ws.Range("B1:B61").NumberFormat
$ws.Columns.Item(3).AutoFit() | Out-Null

$ws.Range("A26").Select() | Out-Null
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Range("C2:C61").Select()
